# Auto-generated edit script
# Updates market-price-derived columns (H, I, J, K, L, M, N) on several
# sheets to reflect a refreshed pull from the market-board data source.
# Source data has no formulas; every cell below is a cached literal value.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 308.48
$ws.Range("I53").Value = 224.375
$ws.Range("K53").Value = 224.375
$ws.Range("M53").Value = 412.625

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1026.3636
$ws.Range("I45").Value = 1039.4286
$ws.Range("K45").Value = 1039.4286
$ws.Range("M45").Value = -662.4286
# Row 46
$ws.Range("H46").Value = 6670.25
$ws.Range("J46").Value = 6670.25
$ws.Range("L46").Value = 6670.25
$ws.Range("N46").Value = -7308.25
# Row 61
$ws.Range("H61").Value = 2113.9473
$ws.Range("I61").Value = 1514.7878
$ws.Range("K61").Value = 1514.7878
$ws.Range("M61").Value = -1302.7878
# Row 74
$ws.Range("H74").Value = 3175.3225
$ws.Range("I74").Value = 950.5833
$ws.Range("J74").Value = 10803
$ws.Range("K74").Value = 950.5833
$ws.Range("L74").Value = 10803
$ws.Range("M74").Value = -76.58330000000001
$ws.Range("N74").Value = -12551
# Row 77
$ws.Range("H77").Value = 3175.3225
$ws.Range("I77").Value = 950.5833
$ws.Range("J77").Value = 10803
$ws.Range("K77").Value = 4752.9165
$ws.Range("L77").Value = 54015
$ws.Range("M77").Value = -384.9165000000003
$ws.Range("N77").Value = -62751
# Row 135
$ws.Range("H135").Value = 33185.8
$ws.Range("J135").Value = 33185.8
$ws.Range("L135").Value = 33185.8
$ws.Range("N135").Value = -43325.8
# Row 136
$ws.Range("H136").Value = 2113.9473
$ws.Range("I136").Value = 1514.7878
$ws.Range("K136").Value = 4544.3634
$ws.Range("M136").Value = -1994.3634

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2062.9644
$ws.Range("I31").Value = 1192.579
$ws.Range("K31").Value = 1192.579
$ws.Range("M31").Value = -897.579
# Row 34
$ws.Range("H34").Value = 2062.9644
$ws.Range("I34").Value = 1192.579
$ws.Range("K34").Value = 1192.579
$ws.Range("M34").Value = -990.579
# Row 99
$ws.Range("H99").Value = 8929604
$ws.Range("I99").Value = 15625828
$ws.Range("K99").Value = 15625828
$ws.Range("M99").Value = -15624330
# Row 126
$ws.Range("H126").Value = 8929604
$ws.Range("I126").Value = 15625828
$ws.Range("K126").Value = 46877484
$ws.Range("M126").Value = -46875014
# Row 134
$ws.Range("H134").Value = 2024.3334
$ws.Range("I134").Value = 1402.5952
$ws.Range("K134").Value = 4207.7856
$ws.Range("M134").Value = -1672.7856

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 2186.182
$ws.Range("I3").Value = 2186.182
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6558.545999999999
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -6446.545999999999
$ws.Range("N3").ClearContents()
# Row 133
$ws.Range("H133").Value = 28336
$ws.Range("I133").Value = 3226.6667
$ws.Range("J133").Value = 66000
$ws.Range("K133").Value = 9680.000100000001
$ws.Range("L133").Value = 198000
$ws.Range("M133").Value = -4620.000100000001
$ws.Range("N133").Value = -208120
# Row 134
$ws.Range("H134").Value = 4817.8184
$ws.Range("I134").Value = 3543.389
$ws.Range("J134").Value = 10552.75
$ws.Range("K134").Value = 10630.167
$ws.Range("L134").Value = 31658.25
$ws.Range("M134").Value = -5560.167000000001
$ws.Range("N134").Value = -41798.25
# Row 136
$ws.Range("H136").Value = 2521.5881
$ws.Range("I136").Value = 1935
$ws.Range("J136").Value = 2841.5454
$ws.Range("K136").Value = 5805
$ws.Range("L136").Value = 8524.636200000001
$ws.Range("M136").Value = -705
$ws.Range("N136").Value = -18724.6362
# Row 138
$ws.Range("H138").Value = 1960
$ws.Range("I138").Value = 786.6667
$ws.Range("J138").Value = 3133.3333
$ws.Range("K138").Value = 2360.0001
$ws.Range("L138").Value = 9399.999899999999
$ws.Range("M138").Value = 2779.9999
$ws.Range("N138").Value = -19679.9999
# Row 139
$ws.Range("H139").Value = 2827.697
$ws.Range("I139").Value = 2605.0435
$ws.Range("J139").Value = 3339.8
$ws.Range("K139").Value = 7815.130500000001
$ws.Range("L139").Value = 10019.4
$ws.Range("M139").Value = -2675.130500000001
$ws.Range("N139").Value = -20299.4
# Row 140
$ws.Range("H140").Value = 4153.5
$ws.Range("I140").Value = 5434.1304
$ws.Range("J140").Value = 2603.2632
$ws.Range("K140").Value = 16302.3912
$ws.Range("L140").Value = 7809.7896
$ws.Range("M140").Value = -11122.3912
$ws.Range("N140").Value = -18169.7896

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 59
$ws.Range("H59").Value = 8363.637000000001
$ws.Range("I59").Value = 7000
$ws.Range("K59").Value = 7000
$ws.Range("M59").Value = -6417
# Row 102
$ws.Range("H102").Value = 2776.4
$ws.Range("I102").Value = 2604
$ws.Range("K102").Value = 2604
$ws.Range("M102").Value = -982
# Row 122
$ws.Range("H122").Value = 1011438.6
$ws.Range("I122").Value = 1853087.5
$ws.Range("J122").Value = 1460
$ws.Range("K122").Value = 5559262.5
$ws.Range("L122").Value = 4380
$ws.Range("M122").Value = -5556812.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2566.2593
$ws.Range("I7").Value = 1537.6154
$ws.Range("K7").Value = 1537.6154
$ws.Range("M7").Value = -1425.6154
# Row 22
$ws.Range("H22").Value = 13736.5
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 18098.666
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 18098.666
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -18688.666
# Row 27
$ws.Range("H27").Value = 13736.5
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 18098.666
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 18098.666
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -18312.666
# Row 40
$ws.Range("H40").Value = 3774
$ws.Range("I40").Value = 1901
$ws.Range("K40").Value = 1901
$ws.Range("M40").Value = -1765
# Row 46
$ws.Range("H46").Value = 1490.7693
$ws.Range("I46").Value = 1182.3529
$ws.Range("J46").Value = 2073.3333
$ws.Range("K46").Value = 1182.3529
$ws.Range("L46").Value = 2073.3333
$ws.Range("M46").Value = -994.3529000000001
$ws.Range("N46").Value = -2449.3333
# Row 55
$ws.Range("H55").Value = 299.79166
$ws.Range("I55").Value = 222
$ws.Range("J55").Value = 355.35715
$ws.Range("K55").Value = 222
$ws.Range("L55").Value = 355.35715
$ws.Range("M55").Value = -49
$ws.Range("N55").Value = -701.35715
# Row 82
$ws.Range("H82").Value = 1325.8
$ws.Range("I82").Value = 1075
$ws.Range("J82").Value = 1493
$ws.Range("K82").Value = 1075
$ws.Range("L82").Value = 1493
$ws.Range("M82").Value = -714
$ws.Range("N82").Value = -2215
# Row 85
$ws.Range("H85").Value = 1325.8
$ws.Range("I85").Value = 1075
$ws.Range("J85").Value = 1493
$ws.Range("K85").Value = 1075
$ws.Range("L85").Value = 1493
$ws.Range("M85").Value = 173
$ws.Range("N85").Value = -3989
# Row 126
$ws.Range("H126").Value = 2566.2593
$ws.Range("I126").Value = 1537.6154
$ws.Range("K126").Value = 4612.8462
$ws.Range("M126").Value = -2142.8462
# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").Value = 3522.4854
$ws.Range("I132").Value = 3584.0784
$ws.Range("K132").Value = 10752.2352
$ws.Range("M132").Value = -8222.235199999999
# Row 136
$ws.Range("H136").Value = 4545.9
$ws.Range("I136").Value = 2700.3225
$ws.Range("J136").Value = 10902.889
$ws.Range("K136").Value = 8100.967500000001
$ws.Range("L136").Value = 32708.667
$ws.Range("M136").Value = -5550.967500000001
$ws.Range("N136").Value = -37808.667

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 145286.28
$ws.Range("I122").Value = 252001
$ws.Range("K122").Value = 756003
$ws.Range("M122").Value = -753553
# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
